$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the "first row" of each 3-row class block gets fresh randomized
# id values; the remaining rows keep their existing text (already correct).
# Written in the same order the shared-string table picks up new entries.
$ws.Range("C5").Value = "EvZxF25Z6fboXyQkPtWX"
$ws.Range("D5").Value = "UhRRtFco6g1LYSAMOF3N"
$ws.Range("E5").Value = "Lcj53mN8KM7991qGd8zz"
$ws.Range("F5").Value = "4A4RT6XxbKqsBEwvV8wo"

$ws.Range("C11").Value = "ifjiFoIfKQhDFiDncZuZ"
$ws.Range("D11").Value = "9a8jgHWcTfM3Z6B12XdT"
$ws.Range("E11").Value = "kCgbyg76CPP5MCNb1peJ"
$ws.Range("F11").Value = "cnQAjUroHCev0i87A4eb"

$ws.Range("C17").Value = "0tR79XBdyBmDnokTgOVL"
$ws.Range("D17").Value = "0UCt9AbEjb0pWbck4ygC"
$ws.Range("E17").Value = "MPCNsmnk79aRe2gZESQo"
$ws.Range("F17").Value = "RDETDHrjCOI2BsNX2709"

$ws.Range("C2").Value = "IW9ZaJ8tACfqeK5l9AxX"
$ws.Range("D2").Value = "pR9LNRfeb2i3sNGCWd5g"
$ws.Range("E2").Value = "oDhHuGW4KUXt9yuuN2sm"
$ws.Range("F2").Value = "YajA9xzXI0UGsaNeTutk"

$ws.Range("C14").Value = "tCXKgjC4Lctq0g9rMCcu"
$ws.Range("D14").Value = "LiyBjLC6w2h4d0ABHT5L"
$ws.Range("E14").Value = "iqx6GLPMejfMddzr82QB"
$ws.Range("F14").Value = "7PGNPORYDCSgawF58iL8"

$ws.Range("C8").Value = "l4cn6ZSRo9yFJSCCkqwh"
$ws.Range("D8").Value = "a56RxfaCjCmUAOYbj4jj"
$ws.Range("E8").Value = "ynI2PlnRWA5r2swEXonq"
$ws.Range("F8").Value = "UWonFliRsXXBsP8vNEkK"

# Move the active selection to F8 (matches the saved view state)
$ws.Range("F8").Select()
